$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BF column ("Date") holds the game date as text. The values were
# originally written as "6-8-2007-08" (month-day-season, ambiguous/incorrect)
# and need to be corrected to the proper ISO-like "2008-06-08" form while
# staying plain text (not an Excel date serial).
$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 58).Value = "2008-06-08"
}
